# Reorder the roster rows (A2:C19) into their new order.
# The player names, positions and team names themselves are unchanged;
# only the order of the rows has been rearranged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Scoot Henderson", "PG", "Portland Trail Blazers"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Walker Kessler", "C", "Utah Jazz")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
